# Display total in overview page
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the Expiry Date value (and its date formatting) for rows 13 and 15 (column F)
$ws.Range("F13").Clear()
$ws.Range("F15").Clear()

# Rename "Tablet" -> "Tablets" in column H for the data rows (2-18)
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    if ($cell.Text -eq "Tablet") {
        $cell.Value = "Tablets"
    }
}

# Delete row 19 entirely (shifts dimension from A1:H19 to A1:H18)
$ws.Rows.Item(19).Delete()
